$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 39
$ws1.Range("G2").Value = 100
$ws1.Range("H2").Value = 7.7

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 6
$ws1.Range("F3").Value = 30
$ws1.Range("G3").Value = 83.33
$ws1.Range("H3").Value = 8.300000000000001

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 7
$ws1.Range("F4").Value = 14
$ws1.Range("G4").Value = 66.67
$ws1.Range("H4").Value = 7

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 12
$ws2.Range("F2").Value = 27
$ws2.Range("G2").Value = 69.23
$ws2.Range("H2").Value = 7.7

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 5
$ws2.Range("F3").Value = 31
$ws2.Range("G3").Value = 86.11
$ws2.Range("H3").Value = 8.300000000000001

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 6
$ws2.Range("F4").Value = 15
$ws2.Range("G4").Value = 71.43000000000001
$ws2.Range("H4").Value = 7

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 12
$ws3.Range("F2").Value = 27
$ws3.Range("G2").Value = 69.23
$ws3.Range("H2").Value = 7.3

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 5
$ws3.Range("F3").Value = 31
$ws3.Range("G3").Value = 86.11
$ws3.Range("H3").Value = 8.4

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 6
$ws3.Range("F4").Value = 15
$ws3.Range("G4").Value = 71.43000000000001
$ws3.Range("H4").Value = 7

# --- Sheet "Rescatables" ---
# Rows 3 and 4 (the students CANO/ORTEGA/OMAR and OCHOA/REYES/OSCAR URIEL) are
# no longer rescatable, so they are removed; the remaining row's "Reprobadas"
# count (G2) drops from 2 to 1.
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Range("G2").Value = 1
$ws4.Rows.Item(3).Resize(2, 1).EntireRow.Delete()
